# Re-generate the statistics with fixed minutes/seconds formatting
# in the "Общее время" (total time) column: zero-pad minutes and
# seconds to two digits (hours stay unpadded), e.g.
#   "1 ч. 9 мин. 36 сек."  -> "1 ч. 09 мин. 36 сек."
#   "3 ч. 50 мин. 6 сек."  -> "3 ч. 50 мин. 06 сек."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the "Общее время" column by inspecting the header row.
$lastCol = $ws.Cells.Item(1, 1).End(2).Column
$timeCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Text
    if ($header -eq "Общее время") {
        $timeCol = $c
        break
    }
}

if ($timeCol -eq 0) {
    $timeCol = 9
}

$lastRow = $ws.Cells.Item(1, $timeCol).End(4).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $timeCol)
    $v = $cell.Text
    if ($v -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
        $h = $Matches[1]
        $mi = $Matches[2].PadLeft(2, '0')
        $se = $Matches[3].PadLeft(2, '0')
        $newv = "$h ч. $mi мин. $se сек."
        if ($newv -ne $v) {
            $cell.Value = $newv
        }
    }
}

Write-Output "Done fixing time formatting up to row $lastRow, column $timeCol"
